$d = $word.ActiveDocument

# Remove the "Dessutom packeteras ..." paragraph from the "Design val"
# section, along with the blank paragraph immediately before it and the
# blank paragraph immediately after it (three paragraphs total), right
# after the "Model-View-Controller (MVC)" paragraph.
$target = "Dessutom packeteras alla klasser i projektet i subpaket till dess relevanta paket för att lätt dela in projektet i olika lätt överskådliga delar."

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Determine the 1-based paragraph index of the matched text so we can
    # grab its immediate neighbours from the document's Paragraphs collection.
    $precedingRange = $d.Range(0, $findRange.Start)
    $paraIndex = $precedingRange.Paragraphs.Count + 1

    $paragraphs = $d.Paragraphs
    $startPara = $paragraphs.Item($paraIndex - 1)
    $endPara = $paragraphs.Item($paraIndex + 1)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
